$wb = $excel.ActiveWorkbook

# Data refresh for 2024-11-06: updates YTD violent-crime counts across
# Citywide Totals, By Neighborhood, and affected individual neighborhood sheets.

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("D2").Value = 86
$ws.Range("E2").Value = 66
$ws.Range("G2").Value = 81
$ws.Range("H2").Value = 97
$ws.Range("D3").Value = 122
$ws.Range("F3").Value = 120
$ws.Range("H3").Value = 136
$ws.Range("C6").Value = 444
$ws.Range("D6").Value = 375
$ws.Range("E6").Value = 417
$ws.Range("F6").Value = 471
$ws.Range("H6").Value = 414
$ws.Range("I6").Value = 465
$ws.Range("J6").Value = 380
$ws.Range("K6").Value = 455
$ws.Range("C7").Value = 591
$ws.Range("D7").Value = 589
$ws.Range("E7").Value = 622
$ws.Range("F7").Value = 676
$ws.Range("G7").Value = 628
$ws.Range("H7").Value = 661
$ws.Range("I7").Value = 779
$ws.Range("J7").Value = 712
$ws.Range("K7").Value = 808

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("F6").Value = 19
$ws.Range("H6").Value = 21
$ws.Range("F7").Value = 44
$ws.Range("H7").Value = 33

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("D2").Value = 1
$ws.Range("D6").Value = 11

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("F4").Value = 3
$ws.Range("F5").Value = 4

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J6").Value = 12
$ws.Range("J7").Value = 27

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("D3").Value = 4
$ws.Range("J4").Value = 12
$ws.Range("D5").Value = 9
$ws.Range("J5").Value = 15

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("H7").Value = 7
$ws.Range("F8").Value = 42
$ws.Range("I14").Value = 2
$ws.Range("D19").Value = 26
$ws.Range("F20").Value = 7
$ws.Range("H23").Value = 6
$ws.Range("J28").Value = 27
$ws.Range("D30").Value = 10
$ws.Range("F36").Value = 44
$ws.Range("H36").Value = 33
$ws.Range("G45").Value = 6
$ws.Range("D50").Value = 11
$ws.Range("C53").Value = 52
$ws.Range("D53").Value = 66
$ws.Range("F53").Value = 70
$ws.Range("H53").Value = 87
$ws.Range("I53").Value = 118
$ws.Range("J54").Value = 8
$ws.Range("K61").Value = 3
$ws.Range("E70").Value = 17
$ws.Range("D82").Value = 9
$ws.Range("J82").Value = 15
$ws.Range("F88").Value = 4
$ws.Range("E96").Value = 10
$ws.Range("C98").Value = 591
$ws.Range("D98").Value = 589
$ws.Range("E98").Value = 622
$ws.Range("F98").Value = 676
$ws.Range("G98").Value = 628
$ws.Range("H98").Value = 661
$ws.Range("I98").Value = 779
$ws.Range("J98").Value = 712
$ws.Range("K98").Value = 808

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("E2").Value = 2
$ws.Range("E6").Value = 10

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("D2").Value = 10
$ws.Range("H2").Value = 11
$ws.Range("C6").Value = 34
$ws.Range("F6").Value = 51
$ws.Range("I6").Value = 75
$ws.Range("C7").Value = 52
$ws.Range("D7").Value = 66
$ws.Range("F7").Value = 70
$ws.Range("H7").Value = 87
$ws.Range("I7").Value = 118

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("H3").Value = 1
$ws.Range("H7").Value = 6

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("D3").Value = 7
$ws.Range("D7").Value = 26

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("H3").Value = 3
$ws.Range("H6").Value = 7

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("D5").Value = 8
$ws.Range("D6").Value = 10

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("J4").Value = 6
$ws.Range("J5").Value = 8

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("G2").Value = 2
$ws.Range("G6").Value = 6

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("E5").Value = 15
$ws.Range("E6").Value = 17

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("F6").Value = 30
$ws.Range("F7").Value = 42

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("H5").Value = 2
$ws.Range("H6").Value = 2

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("F3").Value = 4
$ws.Range("F6").Value = 7
